$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Refresh the "Förändrad" (changed) date column C for all data rows (2-11):
#    46065 -> 46066 (2026-02-12 -> 2026-02-13)
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = 46066
}

# 2) Rows 7-11 got re-sorted (Beteckning/Datum/Area shifted up one slot, with the
#    former row 7 entry wrapping around to the bottom). Capture the original
#    A/B/G values first, then write them back in the new order.
$origA7 = $ws.Cells.Item(7, 1).Value2
$origB7 = $ws.Cells.Item(7, 2).Value2
$origG7 = $ws.Cells.Item(7, 7).Value2

$origA8 = $ws.Cells.Item(8, 1).Value2
$origB8 = $ws.Cells.Item(8, 2).Value2
$origG8 = $ws.Cells.Item(8, 7).Value2

$origA9 = $ws.Cells.Item(9, 1).Value2
$origB9 = $ws.Cells.Item(9, 2).Value2
$origG9 = $ws.Cells.Item(9, 7).Value2

$origA10 = $ws.Cells.Item(10, 1).Value2
$origB10 = $ws.Cells.Item(10, 2).Value2
$origG10 = $ws.Cells.Item(10, 7).Value2

$origA11 = $ws.Cells.Item(11, 1).Value2
$origB11 = $ws.Cells.Item(11, 2).Value2
$origG11 = $ws.Cells.Item(11, 7).Value2

# New row 7 <- old row 9 (A 33037-2025)
$ws.Cells.Item(7, 1).Value = $origA9
$ws.Cells.Item(7, 2).Value = $origB9
$ws.Cells.Item(7, 7).Value = $origG9

# New row 8 <- old row 10 (A 33033-2025)
$ws.Cells.Item(8, 1).Value = $origA10
$ws.Cells.Item(8, 2).Value = $origB10
$ws.Cells.Item(8, 7).Value = $origG10

# New row 9 <- old row 8 (A 6314-2022)
$ws.Cells.Item(9, 1).Value = $origA8
$ws.Cells.Item(9, 2).Value = $origB8
$ws.Cells.Item(9, 7).Value = $origG8

# New row 10 <- old row 11 (A 25610-2024)
$ws.Cells.Item(10, 1).Value = $origA11
$ws.Cells.Item(10, 2).Value = $origB11
$ws.Cells.Item(10, 7).Value = $origG11

# New row 11 <- old row 7 (A 57810-2022)
$ws.Cells.Item(11, 1).Value = $origA7
$ws.Cells.Item(11, 2).Value = $origB7
$ws.Cells.Item(11, 7).Value = $origG7

Write-Host "edit complete"
